# Insert a new localization row ("strChkAbsoluteIntegral") into the
# sorted translation table on the "fr-FR" sheet, just above the existing
# "strChkComputeDerivative" row (row 9), pushing everything else down by
# one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 9; this shifts rows 9..167 down to 10..168
# and inherits formatting from the surrounding rows.
$ws.Rows.Item(9).Insert()

# Populate the new row's data.
$ws.Range("B9").Value = "strChkAbsoluteIntegral"
$ws.Range("C9").Value = 'In "settings" form, tab "Integration"'
$ws.Range("D9").Value = "Compute the absolute-value integral?"

# Grow the table (ListObject) so it includes the newly inserted row.
$tbl = $ws.ListObjects.Item(1)
$newLastRow = $tbl.Range.Row + $tbl.Range.Rows.Count    # one more than before, since a row was inserted inside the table
$tbl.Resize($ws.Range("B2:E" + $newLastRow))
